$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 13.739149
$ws.Range("H2").Value = 41.217447
$ws.Range("I2").Value = 0.6130043224686931
$ws.Range("J2").Value = 0.6130043224686931
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.08532
$ws.Range("N2").Value = 6.25596
$ws.Range("O2").Value = 0.01753772176136817
$ws.Range("P2").Value = 0.01753772176136816
$ws.Range("Q2").Value = 28.65052219268
$ws.Range("R2").Value = 257.85469973412
$ws.Range("S2").Value = 0.01075069924597195
$ws.Range("T2").Value = 0.01075069924597194
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 13.739149
$ws.Range("H3").Value = 41.217447
$ws.Range("I3").Value = 0.6130043224686931
$ws.Range("J3").Value = 0.6130043224686931
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 101.898173
$ws.Range("N3").Value = 305.694519
$ws.Range("O3").Value = 0.8569724579756384
$ws.Range("P3").Value = 0.8569724579756383
$ws.Range("Q3").Value = 1399.994181674777
$ws.Range("R3").Value = 12599.94763507299
$ws.Range("S3").Value = 0.5253278209756868
$ws.Range("T3").Value = 0.5253278209756866
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 13.739149
$ws.Range("H4").Value = 41.217447
$ws.Range("I4").Value = 0.6130043224686931
$ws.Range("J4").Value = 0.6130043224686931
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 14.921347
$ws.Range("N4").Value = 44.76404100000001
$ws.Range("O4").Value = 0.1254898202629935
$ws.Range("P4").Value = 0.1254898202629935
$ws.Range("Q4").Value = 205.006609713703
$ws.Range("R4").Value = 1845.059487423327
$ws.Range("S4").Value = 0.07692580224703441
$ws.Range("T4").Value = 0.07692580224703439
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 4.034036666666666
$ws.Range("H5").Value = 12.10211
$ws.Range("I5").Value = 0.1799879973398545
$ws.Range("J5").Value = 0.1799879973398545
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 2.08532
$ws.Range("N5").Value = 6.25596
$ws.Range("O5").Value = 0.01753772176136817
$ws.Range("P5").Value = 0.01753772176136816
$ws.Range("Q5").Value = 8.412257341733332
$ws.Range("R5").Value = 75.71031607559999
$ws.Range("S5").Value = 0.003156579417732242
$ws.Range("T5").Value = 0.003156579417732241
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 4.034036666666666
$ws.Range("H6").Value = 12.10211
$ws.Range("I6").Value = 0.1799879973398545
$ws.Range("J6").Value = 0.1799879973398545
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 101.898173
$ws.Range("N6").Value = 305.694519
$ws.Range("O6").Value = 0.8569724579756384
$ws.Range("P6").Value = 0.8569724579756383
$ws.Range("Q6").Value = 411.0609661483433
$ws.Range("R6").Value = 3699.54869533509
$ws.Range("S6").Value = 0.1542447564864478
$ws.Range("T6").Value = 0.1542447564864478
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 4.034036666666666
$ws.Range("H7").Value = 12.10211
$ws.Range("I7").Value = 0.1799879973398545
$ws.Range("J7").Value = 0.1799879973398545
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 14.921347
$ws.Range("N7").Value = 44.76404100000001
$ws.Range("O7").Value = 0.1254898202629935
$ws.Range("P7").Value = 0.1254898202629935
$ws.Range("Q7").Value = 60.19326091405667
$ws.Range("R7").Value = 541.7393482265101
$ws.Range("S7").Value = 0.0225866614356745
$ws.Range("T7").Value = 0.02258666143567449
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 4.639623666666666
$ws.Range("H8").Value = 13.918871
$ws.Range("I8").Value = 0.2070076801914524
$ws.Range("J8").Value = 0.2070076801914524
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 2.08532
$ws.Range("N8").Value = 6.25596
$ws.Range("O8").Value = 0.01753772176136817
$ws.Range("P8").Value = 0.01753772176136816
$ws.Range("Q8").Value = 9.675100024573332
$ws.Range("R8").Value = 87.07590022116
$ws.Range("S8").Value = 0.003630443097663977
$ws.Range("T8").Value = 0.003630443097663976
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 4.639623666666666
$ws.Range("H9").Value = 13.918871
$ws.Range("I9").Value = 0.2070076801914524
$ws.Range("J9").Value = 0.2070076801914524
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 101.898173
$ws.Range("N9").Value = 305.694519
$ws.Range("O9").Value = 0.8569724579756384
$ws.Range("P9").Value = 0.8569724579756383
$ws.Range("Q9").Value = 472.7691750408943
$ws.Range("R9").Value = 4254.922575368049
$ws.Range("S9").Value = 0.1773998805135038
$ws.Range("T9").Value = 0.1773998805135038
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 4.639623666666666
$ws.Range("H10").Value = 13.918871
$ws.Range("I10").Value = 0.2070076801914524
$ws.Range("J10").Value = 0.2070076801914524
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 14.921347
$ws.Range("N10").Value = 44.76404100000001
$ws.Range("O10").Value = 0.1254898202629935
$ws.Range("P10").Value = 0.1254898202629935
$ws.Range("Q10").Value = 69.22943467974567
$ws.Range("R10").Value = 623.0649121177111
$ws.Range("S10").Value = 0.02597735658028461
$ws.Range("T10").Value = 0.0259773565802846
